## Menu planning model - add more analyses functions and series run
## Adds three new parametric-sweep worksheets (Stepwise_reduction_carbon,
## Household_size, All_objs) after the existing Stepwise_reduction_waste
## sheet, plus minor selection/view-state tweaks on the first two sheets.

$wb = $excel.ActiveWorkbook
$ws1 = $wb.Worksheets.Item("Sheet1")
$ws2 = $wb.Worksheets.Item("Stepwise_reduction_waste")

function Set-HeaderRow($ws) {
    $ws.Cells.Item(1, 1).Value = "run"
    $ws.Cells.Item(1, 2).Value = "n_days"
    $ws.Cells.Item(1, 3).Value = "n_persons"
    $ws.Cells.Item(1, 4).Value = "dev"
    $ws.Cells.Item(1, 5).Value = "optimize_over"
    $ws.Cells.Item(1, 6).Value = "DRVs"
    $ws.Cells.Item(1, 7).Value = "tvar1"
    $ws.Range("A1:G1").Font.Bold = $true
}

## ---------------------------------------------------------------------
## 1. New sheet: Stepwise_reduction_carbon (rows 1-17, cols A-G)
## ---------------------------------------------------------------------
$ws3 = $wb.Worksheets.Add($null, $ws2)
$ws3.Name = "Stepwise_reduction_carbon"

Set-HeaderRow $ws3

$runNames3 = @("run 1","run 2","run 3","run 4","run 5","run 6","run 7","run 8","run 9","run 10","run 11","run 12","run 13","run 14","run 15","run 16")
$landuse3 = @(29000,28500,28000,27500,27000,26500,26000,25500,25000,24500,24000,23500,23000,22500,22000,21500)

for ($i = 0; $i -lt $runNames3.Count; $i++) {
    $r = $i + 2
    $ws3.Cells.Item($r, 1).Value = $runNames3[$i]
    $ws3.Cells.Item($r, 2).Value = 5
    $ws3.Cells.Item($r, 3).Value = 4
    $ws3.Cells.Item($r, 4).Value = 0.1
    $ws3.Cells.Item($r, 5).Value = "Waste_grams"
    $ws3.Cells.Item($r, 6).Value = "modelgezin_gemiddeld"
    $ws3.Cells.Item($r, 7).Value = $landuse3[$i]
}

$ws3.Columns.Item(3).ColumnWidth = 9.71
$ws3.Columns.Item(5).ColumnWidth = 13.14
$ws3.Columns.Item(6).ColumnWidth = 19.71

$ws3.Range("E2").Select() | Out-Null

## ---------------------------------------------------------------------
## 2. New sheet: Household_size (rows 1-11, cols A-G)
## ---------------------------------------------------------------------
$ws4 = $wb.Worksheets.Add($null, $ws3)
$ws4.Name = "Household_size"

Set-HeaderRow $ws4

$runNames4 = @("run 1","run 2","run 3","run 4","run 5","run 6","run 7","run 8","run 9","run 10")

for ($i = 0; $i -lt $runNames4.Count; $i++) {
    $r = $i + 2
    $ws4.Cells.Item($r, 1).Value = $runNames4[$i]
    $ws4.Cells.Item($r, 2).Value = 5
    $ws4.Cells.Item($r, 3).Value = $i + 1
    $ws4.Cells.Item($r, 4).Value = 0.1
    $ws4.Cells.Item($r, 5).Value = "Waste_grams"
    $ws4.Cells.Item($r, 6).Value = "modelgezin_gemiddeld"
    $ws4.Cells.Item($r, 7).Value = 9999
}

$ws4.Columns.Item(3).ColumnWidth = 10.14
$ws4.Columns.Item(5).ColumnWidth = 14
$ws4.Columns.Item(6).ColumnWidth = 22.43

$ws4.Range("B4:B11").Select() | Out-Null

## ---------------------------------------------------------------------
## 3. New sheet: All_objs (rows 1-6, cols A-G)
## ---------------------------------------------------------------------
$ws5 = $wb.Worksheets.Add($null, $ws4)
$ws5.Name = "All_objs"

Set-HeaderRow $ws5

$runNames5 = @("run 1","run 2","run 3","run 4","run 5")
$optimizeOver5 = @("Waste_grams","Total_cost","Total_carbon","Carbon_waste","Total_landuse")

for ($i = 0; $i -lt $runNames5.Count; $i++) {
    $r = $i + 2
    $ws5.Cells.Item($r, 1).Value = $runNames5[$i]
    $ws5.Cells.Item($r, 2).Value = 5
    $ws5.Cells.Item($r, 3).Value = 4
    $ws5.Cells.Item($r, 4).Value = 0.1
    $ws5.Cells.Item($r, 5).Value = $optimizeOver5[$i]
    $ws5.Cells.Item($r, 6).Value = "modelgezin_gemiddeld"
    $ws5.Cells.Item($r, 7).Value = 9999
}

$ws5.Columns.Item(3).ColumnWidth = 10.14
$ws5.Columns.Item(5).ColumnWidth = 14
$ws5.Columns.Item(6).ColumnWidth = 22.43

$ws5.Range("O8").Select() | Out-Null
$ws5.Activate()

## ---------------------------------------------------------------------
## 4. View-state tweaks on the pre-existing sheets
## ---------------------------------------------------------------------
$ws1.Activate()
$ws1.Range("E2:E5").Select() | Out-Null

$ws2.Activate()
$ws2.Range("A1:G6").Select() | Out-Null

$ws5.Activate()
